$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "meetup" card (column F) : add a new "h3: Meetup coming in" line.
# Written before the "signin" card below so the shared-string table keeps the
# same append order the source workbook ended up with.
$meetupText = @"
type: meetup
width: 2
height: 1
h3: Meetup coming in
date: 1 Mar 2020
button.default: Speak
button.default: Attend
button.default: Details
"@
$ws.Range("F2").Value = $meetupText

# --- Update the "signin" card (column E) : h3.width-half -> h3.w-half.
$signinText = @"
type: signin
width: 2
height: 1
h3.w-half: Sign up to get unlimited access to the entire content of zakatlists
button.primary: Sign In
button.secondary: Sign Up for Rs 300 / Month
"@
$ws.Range("E2").Value = $signinText

# --- View state: scroll the window down so row 2 is at the top, and move the
# active selection to E2 (the signin card).
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E2").Select()
